# Welcome.docx (Chinese Traditional) - apply the authored edit:
#   Trim the trailing clause ", until Smartcash reaches a considerable
#   market cap." from the SmartCash-mining paragraph.

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "and it’s probably safe to assume no ASICs will be created for quite some time, until Smartcash reaches a considerable market cap.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "and it’s probably safe to assume no ASICs will be created for quite some time.",
    2
)
